$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 18: R_ON_ACT label + value
$ws.Range("A18").Value = "R_ON_ACT"
$ws.Range("B18").Value = 44200

# New defined name pointing at the new cell
$wb.Names.Add('R_ON_ACT', '=Sheet1!$B$18')

# Match the number format already used by K (B6) and f_S (B7)
$ws.Range("B18").NumberFormat = $ws.Range("B6").NumberFormat

# New row 19: t_ON label + formula
$ws.Range("A19").Value = "t_ON"
$ws.Range("B19").Formula = "=K*R_ON_ACT/V_IN"

$ws.Range("B19").Select()
